$d = $word.ActiveDocument

# --- Paragraph: "When one shape is selected, you can't select another shape..." ---
# Highlight the entire paragraph (including the paragraph mark) in bright green.
$r1 = $d.Content
$r1.Find.Execute("When one shape is selected, you can", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $r1.Paragraphs(1)
$para1.Range.Font.HighlightColorIndex = 4

# --- Paragraph: "I have to click really far away ... bounding area. (GIF C)" ---
# Highlight everything up to and including "bounding area." in bright green,
# leaving " (GIF C)" without highlight (it becomes its own, unhighlighted run).
$r2 = $d.Content
$r2.Find.Execute("I have to click really far away from the selected object to deselect it. Fixable? At least reduce the size of the bounding area.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Font.HighlightColorIndex = 4
